$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old sample rows (2-6), keep only row 1 ---
$ws.Range("A2:C6").ClearContents()

# --- Row 1 content: new movie row (Avengers / Fantasy / Marvel / date / number) ---
$ws.Range("A1").Value = "Avengers "
$ws.Range("B1").Value = "Fantasy "
$ws.Range("C1").Value = "Marvel "

$ws.Range("D1").NumberFormat = "mm-dd-yy"
$ws.Range("D1").Value = 44107

$ws.Range("E1").NumberFormat = "#,##0"
$ws.Range("E1").Value = 16093323

# --- Column widths (A:D narrower, E:G a touch wider) ---
$ws.Range("A1:D1").ColumnWidth = 15.33
$ws.Range("E1:G1").ColumnWidth = 16

# --- View: zoom + new selection ---
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("H2").Select()

# --- Comment on A1 (author template Excel auto-fills) ---
$comment = $ws.Range("A1").AddComment("Lucian Schipor:" + [char]10)
